# Refactor for recalc with  (WIP)
#
# - Sheet3!B8: SEARCH("A","ABC") -> SEARCH("B","ABC")  (1 -> 2)
# - Active sheet switches from Sheet1 to Sheet3 (tab selection + per-sheet
#   cell selections move accordingly)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Update the formula on Sheet3 so SEARCH looks for "B" instead of "A".
$ws3.Range("B8").Formula = "=SEARCH(""B"",""ABC"")"

# Leave a cell selected on Sheet1 (no longer the active tab).
$ws1.Activate()
$ws1.Range("B4").Select()

# Sheet3 becomes the active/visible tab, with B8 selected.
$ws3.Activate()
$ws3.Range("B8").Select()
